$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the existing date style already applied to A2:A16 (numFmtId 14,
# "mm-dd-yy") for the newly added date cells in column A, instead of minting
# brand-new (duplicate) style records via a fresh NumberFormat assignment.
$ws.Range("A16").Copy()
$ws.Range("A17:A24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The row that used to be A17/B17 (11/20/2025, error count 1) is pushed down
# to row 20 as new data is appended/interleaved.
$ws.Range("A20").Value = 45981
$ws.Range("B20").Value = 1

# Newly appended "Date" / "Error Count" rows.
$ws.Range("A17").Value = 45982
$ws.Range("B17").Value = 6

$ws.Range("A18").Value = 45987
$ws.Range("B18").Value = 3

$ws.Range("A19").Value = 45988
$ws.Range("B19").Value = 1

$ws.Range("A21").Value = 45989
$ws.Range("B21").Value = 1

$ws.Range("A22").Value = 45983
$ws.Range("B22").Value = 2

$ws.Range("A23").Value = 45985
$ws.Range("B23").Value = 4

$ws.Range("A24").Value = 45986
$ws.Range("B24").Value = 0

# Leave the selection the way the author left it when they saved the file.
$ws.Range("B25").Select()
